# Generate Report for Archive
#
# The localization status for e2e\1eadbff6-f6a7-4f04-88ee-ac3d2e4c5549.md
# advanced from "Ready for handoff" to "In Translation". Update the
# Status columns for that file's row (row 5) on all three worksheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns E (zh-cn) and F (de-de) hold the status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"

# --- zh-cn sheet: column C holds the Status ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "In Translation"

# --- de-de sheet: column C holds the Status ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "In Translation"
